$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("lido")

$ws.Range("A37").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.03857596179489347
$ws.Range("E2").Value = 0.0003906250000000888
$ws.Range("D3").Value = 0.02180615831256907
$ws.Range("E3").Value = 0.00196001568012516
$ws.Range("D4").Value = 0.02015594003783184
$ws.Range("E4").Value = -0.001652892561983421
$ws.Range("D5").Value = 0.04077435325400019
$ws.Range("E5").Value = -0.002795248078266921
$ws.Range("D6").Value = 0.03760164478257103
$ws.Range("E6").Value = 0.0007840062720501084
$ws.Range("D7").Value = 0.02107020129145087
$ws.Range("E7").Value = 0.0007801833430856053
$ws.Range("D8").Value = 0.03777814076686272
$ws.Range("E8").Value = -0.004946043165467651
$ws.Range("D9").Value = 0.02160626469235917
$ws.Range("E9").Value = -0.003560668310052084
$ws.Range("D10").Value = 0.02551273836889546
$ws.Range("E10").Value = 0.009664948453608213
$ws.Range("D11").Value = 0.02413737382194684
$ws.Range("E11").Value = -0.01049943246311014
$ws.Range("D12").Value = 0.05881409493314058
$ws.Range("E12").Value = -0.01145522388059717
$ws.Range("D13").Value = 0.02649248012584531
$ws.Range("E13").Value = 0.002233804914370863
$ws.Range("D14").Value = 0.02682004702870024
$ws.Range("E14").Value = 0.01574803149606296
$ws.Range("D15").Value = 0.0347297630586617
$ws.Range("E15").Value = 0.005654709312599371
$ws.Range("D16").Value = 0.01913455377667159
$ws.Range("E16").Value = -0.01777777777777778
$ws.Range("D17").Value = 0.02964157177975515
$ws.Range("E17").Value = 0.006654958868657079
$ws.Range("D18").Value = 0.02407238647721851
$ws.Range("E18").Value = 0.002321801718133232
$ws.Range("D19").Value = 0.1337914082675409
$ws.Range("E19").Value = 0.0006724949562879612
$ws.Range("D20").Value = 0.00976854039355262
$ws.Range("E20").Value = -0.0159979357502259
$ws.Range("D21").Value = 0.01566912827696357
$ws.Range("E21").Value = 0.01358595872806068
$ws.Range("D22").Value = 0.01724526317052652
$ws.Range("E22").Value = -0.002052611678280192
$ws.Range("D23").Value = 0.01704608189052021
$ws.Range("E23").Value = 0.002102312543798224
$ws.Range("D24").Value = 0.02140374089293598
$ws.Range("E24").Value = 0.0001100836635841151
$ws.Range("D25").Value = 0.01207855008302325
$ws.Range("E25").Value = -0.03860182370820653
$ws.Range("D26").Value = 0.04321548833633329
$ws.Range("E26").Value = 0.002282323405226494
$ws.Range("D27").Value = 0.02568983710259019
$ws.Range("E27").Value = 0.0001962323390896081
$ws.Range("D28").Value = 0.04787616590227996
$ws.Range("E28").Value = -0.0012189176011701
$ws.Range("D29").Value = 0.05774843398855663
$ws.Range("E29").Value = -0.0104147294030128
$ws.Range("D30").Value = 0.01367134277827036
$ws.Range("E30").Value = -0.04258150365934787
$ws.Range("D31").Value = 0.01465776957405398
$ws.Range("E31").Value = -0.002803738317757043
$ws.Range("D32").Value = 0.04445627538020199
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 0.01695829965927672
$ws.Range("E33").Value = -0.01456618112729557
$ws.Range("D34").Value = 0.9999999999999999
$ws.Range("E34").Value = -0.002140272717803149

$ws.Protect("lido")